# Update localization status for the "5104585e-3894-4031-bd11-129864ca0782.md" row
# from "Ready for handoff" to "In Translation" on every sheet of the workbook
# (Overview summary columns B4/C4, and the "Status" column (C4) on the
# per-locale "zh-cn" and "de-de" sheets).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 4 corresponds to 5104585e-...md, columns B (zh-cn) and C (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

# --- zh-cn sheet: row 4 corresponds to 5104585e-...md, column C is "Status" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet: row 4 corresponds to 5104585e-...md, column C is "Status" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C4").Value = "In Translation"
